# Add the two new trailing columns (評估淨值 / 貸放成數) to the header row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "評估淨值"
$ws.Range("O1").Value = "貸放成數"

# Zoom the view to 55% and leave the selection on the newly-touched cell,
# matching the author's saved view state.
$excel.ActiveWindow.Zoom = 55
[void]$ws.Range("O2").Select()
